$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.307.38'
$ws.Range('E2').Value = '  +2.85%  '
$ws.Range('D3').Value = '2.303.47'
$ws.Range('E3').Value = '  +2.00%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '''308.75'
$ws.Range('E5').Value = '  +0.37%  '
$ws.Range('D6').Value = '''104.35'
$ws.Range('E6').Value = '  +7.59%  '
$ws.Range('D7').Value = '''0.529'
$ws.Range('E7').Value = '  +0.79%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').Value = '''0.525'
$ws.Range('E9').Value = '  +7.61%  '
$ws.Range('D10').Value = '''36.20'
$ws.Range('E10').Value = '  +4.32%  '
$ws.Range('D11').Value = '''52.62'
$ws.Range('E11').Value = '  +0.93%  '
$ws.Range('D12').Value = '''0.0810'
$ws.Range('E12').Value = '  -1.24%  '
$ws.Range('E13').Value = '  -0.70%  '
$ws.Range('D14').Value = '''6.95'
$ws.Range('E14').Value = '  +2.23%  '
$ws.Range('D15').Value = '2.660.20'
$ws.Range('E15').Value = '  +1.97%  '
$ws.Range('D16').Value = '''15.03'
$ws.Range('E16').Value = '  +3.11%  '
$ws.Range('D17').Value = '2.301.69'
$ws.Range('E17').Value = '  +1.57%  '
$ws.Range('D18').Value = '''0.802'
$ws.Range('E18').Value = '  +2.33%  '
$ws.Range('D19').Value = '43.234.80'
$ws.Range('E19').Value = '  +2.98%  '
$ws.Range('D20').Value = '''11.89'
$ws.Range('E20').Value = '  -2.92%  '
$ws.Range('D21').Value = '0.0₃0922'
$ws.Range('E21').Value = '  +2.10%  '
$ws.Range('D22').Value = '''6.15'
$ws.Range('E22').Value = '  +3.76%  '
$ws.Range('D23').Value = '''67.81'
$ws.Range('E23').Value = '  +0.59%  '
$ws.Range('D24').Value = '''240.47'
$ws.Range('E24').Value = '  +2.09%  '
$ws.Range('D25').Value = '''2.01'
$ws.Range('E25').Value = '  +2.54%  '
$ws.Range('D26').Value = '''2.60'
$ws.Range('E26').Value = '  +0.62%  '
$ws.Range('D27').Value = '''0.999'
$ws.Range('E27').Value = '  -0.08%  '
$ws.Range('D28').Value = '''24.78'
$ws.Range('E28').Value = '  +5.58%  '
$ws.Range('B29').Value = 'InjectiveProtocol'
$ws.Range('C29').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D29').Value = '''36.29'
$ws.Range('E29').Value = '  -1.26%  '
$ws.Range('B30').Value = 'Cosmos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D30').Value = '''9.56'
$ws.Range('E30').Value = '  +0.34%  '
$ws.Range('B31').Value = 'Toncoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D31').Value = '''2.12'
$ws.Range('E31').Value = '  -0.64%  '
$ws.Range('D32').Value = '''163.75'
$ws.Range('E32').Value = '  -0.48%  '
$ws.Range('D33').Value = '''5.22'
$ws.Range('E33').Value = '  +0.43%  '
$ws.Range('E34').Value = '  +0.03%  '
$ws.Range('D35').Value = '''18.24'
$ws.Range('E35').Value = '  +4.42%  '
$ws.Range('E36').Value = '  +7.00%  '
$ws.Range('D37').Value = '''0.0735'
$ws.Range('E37').Value = '  +1.28%  '
$ws.Range('D38').Value = '''3.00'
$ws.Range('E38').Value = '  -2.57%  '
$ws.Range('D39').Value = '''4.48'
$ws.Range('E39').Value = '  +8.33%  '
$ws.Range('D40').Value = '''1.87'
$ws.Range('E40').Value = '  +3.66%  '
$ws.Range('D41').Value = '''0.105'
$ws.Range('E41').Value = '  +1.94%  '
$ws.Range('E42').Value = '  +0.31%  '
$ws.Range('D43').Value = '''2.46'
$ws.Range('E43').Value = '  +11.15%  '
$ws.Range('D44').Value = '1.982.22'
$ws.Range('E44').Value = '  +2.20%  '
$ws.Range('D45').Value = '''0.0289'
$ws.Range('E45').Value = '  +2.65%  '
$ws.Range('D46').Value = '''18.85'
$ws.Range('E46').Value = '  +1.06%  '
$ws.Range('D47').Value = '''3.06'
$ws.Range('E47').Value = '  +4.60%  '
$ws.Range('D48').Value = '''10.07'
$ws.Range('E48').Value = '  +4.04%  '
$ws.Range('D49').Value = '''57.37'
$ws.Range('E49').Value = '  +6.11%  '
$ws.Range('B50').Value = 'HuobiToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D50').Value = '''2.91'
$ws.Range('E50').Value = '  +1.01%  '
$ws.Range('B51').Value = 'Stacks'
$ws.Range('C51').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D51').Value = '''1.58'
$ws.Range('E51').Value = '  +7.98%  '
